# Daily attendance processing - reorders the "Recorded By" (column G) list so
# that the System/system account is no longer always listed first: the first
# two comma-separated entries in each multi-entry cell are swapped. Cells
# whose value is exactly "System, admin@admin.com" are left untouched, as are
# single-entry cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -eq "System, admin@admin.com") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Length -ge 2) {
        $first = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $first
        $cell.Value = $parts -join ", "
    }
}
